# Add "2022-Q4" quarter data to the 000039-中集集团 workbook.
#
# Two changes, mirroring the commit's xml diff:
#   1. A brand-new worksheet "2022-Q4" is inserted right after "总计",
#      pushing every other quarter sheet down one slot.
#   2. The "总计" (summary) sheet gets a new row 2 for the 2022-Q4
#      quarter, and the pre-existing rows shift down by one.
#
# NOTE: worksheet object handles returned by Worksheets.Item(...) are
# positional - inserting a new sheet shifts what an already-held handle
# refers to. So every sheet reference we need *after* inserting "2022-Q4"
# is (re-)fetched by name only after that insert has happened.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTotal)
$newSheet.Name = "2022-Q4"

# Fetch the template sheet (used only to copy cell formatting from) AFTER
# the insert, since it shifted from slot 2 to slot 3.
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# Header row (row 1, columns B:H) - same header labels as every other
# quarter sheet.
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Single data row (row 2) - the one fund held that quarter. Leading
# apostrophes force the numeric-looking codes/percentages to be stored
# as text, matching every other quarter sheet.
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "'005120"
$newSheet.Cells.Item(2, 3).Value = "上投摩根量化多因子灵活配置混合"
$newSheet.Cells.Item(2, 4).Value = "'0.19"
$newSheet.Cells.Item(2, 5).Value = "'94.61"
$newSheet.Cells.Item(2, 6).Value = "'2.50"
$newSheet.Cells.Item(2, 7).Value = "'0.0048"
$newSheet.Cells.Item(2, 8).Value = 2

# Match formatting used by every other quarter sheet: bold/bordered
# header row, and a styled index cell in column A.
$wsQ2.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$wsQ2.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

# ---------------------------------------------------------------------
# 2. Insert a new row 2 in "总计" for the 2022-Q4 totals, shifting the
#    existing quarters down by one row.
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").Style = "Normal"   # row Insert() copies the header's bold style; strip it back off

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 1
$wsTotal.Cells.Item(2, 4).Value = 0

# Give the new A2 the same styled-index look as the rest of column A.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# Column A is a plain sequential index (0,1,2,...), independent of the
# quarter label in column B - renumber it for every shifted-down row.
for ($r = 3; $r -le 9; $r++) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
}
